$d = $word.ActiveDocument

# List of (old text, new text) replacements. A WordOpenXML round-trip per
# paragraph (rather than Find.Execute's Replace, or directly assigning
# Range.Text) is used here: several paragraphs in this document start with
# an empty run (<w:r/>) followed by the run that holds the visible text,
# and a plain text-replace on those paragraphs causes the empty run to be
# silently dropped / merged away by the engine. Reading the paragraph's own
# WordOpenXML, patching the text in that XML string, and feeding it back via
# InsertXML keeps the paragraph's run structure (including empty runs)
# completely intact. InsertXML also stamps the rewritten <w:p> with fresh
# w14:paraId/w14:textId/w:rsidR/w:rsidRDefault attributes, so those are
# stripped back out before they're written to the paragraph.
$replacements = @(
    @("Play Cirque Du Soleil Amaluna Free Slot Game - Review", "Play Cirque Du Soleil Amaluna for Free - Review"),
    @("Visually pleasing graphics and dreamlike atmosphere", "Visually pleasing graphics with a dreamlike atmosphere"),
    @("Players can trigger a bonus round with free spins", "Access to bonus mode with free spins"),
    @("The maximum value can be reached with 11 or more Amaluna symbols", "Unique theme and visuals"),
    @("Symbols lack impact given the theme of Cirque Du Soleil", "Could have more impact in terms of symbols and theme"),
    @("Gameplay mechanics are not as unique as the overall theme", "Gameplay mechanics differ from similar Bally slot"),
    @("Read our review of Cirque Du Soleil Amaluna and play this free slot game. Get free spins and try to match as many Amaluna symbols as possible.", "Read our review of Cirque Du Soleil Amaluna and play for free. Discover its unique theme and gameplay mechanics.")
)

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    foreach ($pair in $replacements) {
        $old = $pair[0]
        $new = $pair[1]
        if ($text -like ("*" + $old + "*")) {
            $xml = $p.Range.WordOpenXML
            $newxml = $xml.Replace($old, $new)
            $newxml = $newxml -replace ' w14:paraId="[0-9A-Fa-f]+" w14:textId="[0-9A-Fa-f]+" w:rsidR="[0-9A-Fa-f]+" w:rsidRDefault="[0-9A-Fa-f]+"', ''
            $p.Range.InsertXML($newxml) | Out-Null
        }
    }
}
